$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark from its old location
#        (paragraph "... numId 7 ... для ввода данных,") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

# --- 2. Insert a manual page break before the "Требования, ..." heading
#        paragraph, and move the _GoBack bookmark + lastRenderedPageBreak
#        marker into that heading paragraph (mirrors what Word itself
#        records when a Ctrl+Enter page break is typed at that spot). ---
$rng = $d.Content
$rng.Find.Execute("Требования, предъявляемые", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1).Range

$oxml = $para.WordOpenXML
$bodyStart = $oxml.IndexOf("<w:body>")
$bodyEnd = $oxml.IndexOf("</w:body>")
$bodyContent = $oxml.Substring($bodyStart + 8, $bodyEnd - $bodyStart - 8)
$m = [regex]::Match($bodyContent, '(?s)^<w:p\b.*?</w:p>')
$paraXml = $m.Value

$pPrEndIdx = $paraXml.IndexOf("</w:pPr>") + 8
$firstTIdx = $paraXml.IndexOf("<w:t>")

$modified = $paraXml.Substring(0, $firstTIdx) + "<w:lastRenderedPageBreak/>" + $paraXml.Substring($firstTIdx)
$modified = $modified.Substring(0, $pPrEndIdx) + '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/>' + $modified.Substring($pPrEndIdx)

$pageBreakPara = '<w:p><w:pPr><w:spacing w:line="312" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/></w:rPr><w:br w:type="page"/></w:r></w:p>'

$fullBody = $pageBreakPara + $modified

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fullBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$collapsed = $para.Duplicate
$collapsed.Collapse(1)
$collapsed.InsertXML($xmlFrag)

# --- 3. Remove the now-stale lastRenderedPageBreak from the
#        "Содержание - " paragraph (that paragraph no longer starts
#        a rendered page once the break above was inserted). ---
$rng2 = $d.Content
$rng2.Find.Execute("Содержание - ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$contentsPara = $rng2.Paragraphs(1).Range
$oxml2 = $contentsPara.WordOpenXML
$bodyStart2 = $oxml2.IndexOf("<w:body>")
$bodyEnd2 = $oxml2.IndexOf("</w:body>")
$bodyContent2 = $oxml2.Substring($bodyStart2 + 8, $bodyEnd2 - $bodyStart2 - 8)
$m2 = [regex]::Match($bodyContent2, '(?s)^<w:p\b.*?</w:p>')
$contentsParaXml = $m2.Value
$contentsParaXml = $contentsParaXml.Replace("<w:lastRenderedPageBreak/>", "")

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $contentsParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$collapsed2 = $contentsPara.Duplicate
$collapsed2.Collapse(1)
$collapsed2.InsertXML($xmlFrag2)

# --- 4. Bump the footer's cached PAGE field result from 1 to 2 ---
$footer = $d.Sections(1).Footers(1)
$frng = $footer.Range.Duplicate
$frng.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)

Write-Output "done"
